# Applies the "Fixed a few typos." edit to bumpMapping.docx
$d = $word.ActiveDocument

# 1) "Nicola " + "Frachesen" (spell-checked run split) + " " -> single run
#    "Nicola Frachesen " (text is unchanged, this just collapses the
#    proofErr-wrapped runs back into one run).
$d.Content.Find.Execute(
    "Nicola Frachesen ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Nicola Frachesen ", 2) | Out-Null

# 2) Citation typo: "[1], [2]" -> "[1] - [5]" in the intro paragraph.
$d.Content.Find.Execute(
    "illusion of depth by altering their appearance [1], [2].",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "illusion of depth by altering their appearance [1] - [5].", 2) | Out-Null

# 3) Rewrite the "Bump Mapping Techniques in depth" paragraph (typo / wording
#    fixes throughout).
$oldTechniques = "Typically, bump maps are grayscale images that are limited to 8-bits of color information. In which is calculated to be only 256 variations of black, gray or white. These information can be important is that when values in a bump map are close to 50% gray, there" + [char]8217 + "s little to no details that comes through on the surface. When values get brighter and working its way toward the white spectrum, details appear to pull out of surface. Applying the same logic when values get darker and closer to black, they appear to be pushing into the surface. Due to these techniques, it is important to note bump mapping works best for creating tiny details on a model, for example: pores or wrinkles on skin, and silhouette of the geometry that the bump map is applied to will always be unaffected by the map. "

$newTechniques = "Typically, bump maps are grayscale images that are limited to 8-bits of color information. Thus, only 256 variations of black, gray or white can be calculated. For example, when values in a bump map are close to 50% gray, there" + [char]8217 + "s little to no detail that comes through on the surface. When values get brighter (closer to pure white), details appear more and the surface seems to pop out. Applying the same logic, when values get darker and closer to black they appear to be pushing into the surface. These techniques mean that bump mapping works best for creating tiny details on a model such as pores or wrinkles on skin. However, with few exceptions, the silhouette of the geometry that the bump map is applied to will be unaffected by it. "

$d.Content.Find.Execute(
    $oldTechniques, $true, $false, $false, $false, $false,
    $true, 1, $false, $newTechniques, 2) | Out-Null

# 4) The "_GoBack" last-edit-position bookmark follows the final edit: move
#    it to just after the rewritten paragraph's last sentence, i.e. right
#    after "...unaffected by it." and before the trailing space.
$marker = "unaffected by it."
$searchRange = $d.Content
$searchRange.Find.Execute($marker) | Out-Null
$pos = $searchRange.End
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
